# Applies the "Added product-pool in cms" commit:
#  - marks one item as done on both the "TODO" and "TODO CMS" sheets
#  - appends two new open TODO items on each of those sheets
#  - inserts a new empty "Testfälle CMS" sheet before "TODO CMS"
#  - renames "Testfälle" to "Testfälle Shop"
#
# Two host quirks drive the ordering below:
#  1. Worksheets.Item(...) hands back a *positional* reference, not a stable
#     object identity - once the sheet collection is reordered/resized, any
#     reference grabbed beforehand can silently start pointing at whatever
#     sheet now sits at that old index. So every cell-level edit happens
#     first, while sheet order/count is untouched, and the sheet
#     insertion/renaming (which shuffles indices) happens last.
#  2. New shared-string entries are appended in first-use order, and the
#     workbook's new strings (indices 90-93) are interleaved across the two
#     sheets, so the four new "A column" values must be written in the exact
#     sequence below to land on the right shared-string index.

$wb = $excel.ActiveWorkbook
$todo = $wb.Worksheets.Item("TODO")
$todoCms = $wb.Worksheets.Item("TODO CMS")

# --- 1. New shared strings, in first-use order ------------------------------
$todoCms.Range("A23").Value = "Datenbank-User einschränken auf Mindestrechte"
$todo.Range("A28").Value = "Datenbank-User beschränken auf Mindestrechte "
$todo.Range("A29").Value = "E-Mail-Texte in Datenbank?"
$todoCms.Range("A24").Value = "E-Mail-Text in Datenbank? "

# --- 2. "TODO" sheet B-column updates ---------------------------------------
# Row 19 flips from "offen" to "done"
$todo.Range("B2").Copy()
$todo.Range("B19").PasteSpecial(-4122)
$todo.Range("B19").Value = "done"

# New rows 28/29 get the "offen" styling
$todo.Range("B27").Copy()
$todo.Range("B28").PasteSpecial(-4122)
$todo.Range("B28").Value = "offen"

$todo.Range("B27").Copy()
$todo.Range("B29").PasteSpecial(-4122)
$todo.Range("B29").Value = "offen"

$todo.Range("C29").Select()

# --- 3. "TODO CMS" sheet B-column updates -----------------------------------
# Row 20 flips from "offen" to "done"
$todoCms.Range("B2").Copy()
$todoCms.Range("B20").PasteSpecial(-4122)
$todoCms.Range("B20").Value = "done"

# New rows 23/24 get the "offen" styling
$todoCms.Range("B21").Copy()
$todoCms.Range("B23").PasteSpecial(-4122)
$todoCms.Range("B23").Value = "offen"

$todoCms.Range("B21").Copy()
$todoCms.Range("B24").PasteSpecial(-4122)
$todoCms.Range("B24").Value = "offen"

$todoCms.Range("B20").Select()

# --- 4. Sheet restructuring -------------------------------------------------
# Make "TODO CMS" the active sheet so Worksheets.Add() (which inserts right
# before the active sheet) drops the new sheet in the right slot.
$wb.Worksheets.Item("TODO CMS").Activate()
$newCmsTests = $wb.Worksheets.Add()
$newCmsTests.Name = "Testfälle CMS"

$wb.Worksheets.Item("Testfälle").Name = "Testfälle Shop"

# Adding a sheet makes it active; restore "TODO CMS" as the active/selected
# tab so activeTab / tabSelected land back where the diff expects.
$wb.Worksheets.Item("TODO CMS").Activate()
